$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 350
$ws.Range("I6").Value = 350
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1050
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -938
$ws.Range("N6").ClearContents()
$ws.Range("H11").Value = 1232.7693
$ws.Range("I11").Value = 1232.7693
$ws.Range("K11").Value = 1232.7693
$ws.Range("M11").Value = -1092.7693
$ws.Range("H15").Value = 2237.3857
$ws.Range("I15").Value = 2237.3857
$ws.Range("K15").Value = 6712.157099999999
$ws.Range("M15").Value = -6543.157099999999
$ws.Range("H17").Value = 3073.9736
$ws.Range("J17").Value = 3228.361
$ws.Range("L17").Value = 9685.082999999999
$ws.Range("N17").Value = -10021.083
$ws.Range("H121").Value = 960
$ws.Range("I121").Value = 600
$ws.Range("J121").Value = 1200
$ws.Range("K121").Value = 1800
$ws.Range("L121").Value = 3600
$ws.Range("M121").Value = -53
$ws.Range("N121").Value = -7094
$ws.Range("H137").Value = 2193.7944
$ws.Range("I137").Value = 1735.551
$ws.Range("J137").Value = 3129.375
$ws.Range("K137").Value = 5206.653
$ws.Range("L137").Value = 9388.125
$ws.Range("M137").Value = -2656.653
$ws.Range("N137").Value = -14488.125
$ws.Range("H138").Value = 2682.3264
$ws.Range("I138").Value = 1329.6333
$ws.Range("J138").Value = 4818.1577
$ws.Range("K138").Value = 3988.8999
$ws.Range("L138").Value = 14454.4731
$ws.Range("M138").Value = 1151.1001
$ws.Range("N138").Value = -24734.4731
$ws.Range("H141").Value = 4243.9585
$ws.Range("I141").Value = 3275
$ws.Range("J141").Value = 5600.5
$ws.Range("K141").Value = 9825
$ws.Range("L141").Value = 16801.5
$ws.Range("M141").Value = -4645
$ws.Range("N141").Value = -27161.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19674.598
$ws.Range("I32").Value = 21272.697
$ws.Range("J32").Value = 11538.818
$ws.Range("K32").Value = 21272.697
$ws.Range("L32").Value = 11538.818
$ws.Range("M32").Value = -20985.697
$ws.Range("N32").Value = -12112.818
$ws.Range("H63").Value = 2688.125
$ws.Range("I63").Value = 2301.6667
$ws.Range("J63").Value = 2920
$ws.Range("K63").Value = 2301.6667
$ws.Range("L63").Value = 2920
$ws.Range("M63").Value = -1615.6667
$ws.Range("N63").Value = -4292
$ws.Range("H66").Value = 2688.125
$ws.Range("I66").Value = 2301.6667
$ws.Range("J66").Value = 2920
$ws.Range("K66").Value = 11508.3335
$ws.Range("L66").Value = 14600
$ws.Range("M66").Value = -8076.333500000001
$ws.Range("N66").Value = -21464
$ws.Range("H68").Value = 38099
$ws.Range("J68").Value = 38099
$ws.Range("L68").Value = 38099
$ws.Range("N68").Value = -39721
$ws.Range("H71").Value = 38099
$ws.Range("J71").Value = 38099
$ws.Range("L71").Value = 114297
$ws.Range("N71").Value = -122409
$ws.Range("H75").Value = 40173
$ws.Range("J75").Value = 40173
$ws.Range("L75").Value = 40173
$ws.Range("N75").Value = -41921
$ws.Range("H78").Value = 40173
$ws.Range("J78").Value = 40173
$ws.Range("L78").Value = 120519
$ws.Range("N78").Value = -129255
$ws.Range("H121").Value = 41900
$ws.Range("J121").Value = 41900
$ws.Range("L121").Value = 41900
$ws.Range("N121").Value = -45394
$ws.Range("H122").Value = 2509.2307
$ws.Range("I122").Value = 2693.3333
$ws.Range("J122").Value = 2351.4285
$ws.Range("K122").Value = 8079.999899999999
$ws.Range("L122").Value = 7054.2855
$ws.Range("M122").Value = -5629.999899999999
$ws.Range("N122").Value = -11954.2855
$ws.Range("H126").Value = 6833.3335
$ws.Range("I126").Value = 6833.3335
$ws.Range("K126").Value = 20500.0005
$ws.Range("M126").Value = -18030.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 68300
$ws.Range("J13").Value = 68300
$ws.Range("L13").Value = 68300
$ws.Range("N13").Value = -68636
$ws.Range("H58").Value = 34840
$ws.Range("J58").Value = 34840
$ws.Range("L58").Value = 34840
$ws.Range("N58").Value = -35428
$ws.Range("H128").Value = 4700
$ws.Range("I128").Value = 4700
$ws.Range("K128").Value = 14100
$ws.Range("M128").Value = -11610

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 63000
$ws.Range("J52").Value = 63000
$ws.Range("L52").Value = 63000
$ws.Range("N52").Value = -63588
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H98").Value = 20000
$ws.Range("I98").Value = 20000
$ws.Range("K98").Value = 20000
$ws.Range("M98").Value = -17754
$ws.Range("H115").Value = 43000
$ws.Range("J115").Value = 43000
$ws.Range("L115").Value = 43000
$ws.Range("N115").Value = -45350
$ws.Range("H124").Value = 40660
$ws.Range("J124").Value = 40660
$ws.Range("L124").Value = 40660
$ws.Range("N124").Value = -45570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 713.10345
$ws.Range("I113").Value = 710.2041
$ws.Range("J113").Value = 728.8889
$ws.Range("K113").Value = 2130.6123
$ws.Range("L113").Value = 2186.6667
$ws.Range("M113").Value = 39.38769999999977
$ws.Range("N113").Value = -6526.6667
$ws.Range("H132").Value = 1609
$ws.Range("J132").Value = 1349.1666
$ws.Range("L132").Value = 12142.4994
$ws.Range("N132").Value = -17202.4994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12760
$ws.Range("I80").Value = 18600
$ws.Range("K80").Value = 18600
$ws.Range("M80").Value = -17602
$ws.Range("H82").Value = 34131
$ws.Range("J82").Value = 34131
$ws.Range("L82").Value = 34131
$ws.Range("N82").Value = -34897
$ws.Range("H83").Value = 12760
$ws.Range("I83").Value = 18600
$ws.Range("K83").Value = 93000
$ws.Range("M83").Value = -88008
$ws.Range("H85").Value = 34131
$ws.Range("J85").Value = 34131
$ws.Range("L85").Value = 34131
$ws.Range("N85").Value = -36783
$ws.Range("H117").Value = 29000
$ws.Range("J117").Value = 29000
$ws.Range("L117").Value = 29000
$ws.Range("N117").Value = -35884

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H63").Value = 44542.5
$ws.Range("J63").Value = 44542.5
$ws.Range("L63").Value = 44542.5
$ws.Range("N63").Value = -46040.5
$ws.Range("H64").Value = 31716.666
$ws.Range("J64").Value = 31716.666
$ws.Range("L64").Value = 31716.666
$ws.Range("N64").Value = -32166.666
$ws.Range("H66").Value = 44542.5
$ws.Range("J66").Value = 44542.5
$ws.Range("L66").Value = 133627.5
$ws.Range("N66").Value = -141115.5
$ws.Range("H67").Value = 31716.666
$ws.Range("J67").Value = 31716.666
$ws.Range("L67").Value = 31716.666
$ws.Range("N67").Value = -33276.666
$ws.Range("H117").Value = 61500
$ws.Range("J117").Value = 61500
$ws.Range("L117").Value = 61500
$ws.Range("N117").Value = -70678

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2288.1943
$ws.Range("I132").Value = 1152.238
$ws.Range("J132").Value = 3878.5334
$ws.Range("K132").Value = 3456.714
$ws.Range("L132").Value = 11635.6002
$ws.Range("M132").Value = -926.7139999999999
$ws.Range("N132").Value = -16695.6002
$ws.Range("H137").Value = 30000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
